# Update gh-pages output data (苏州-漫展信息.xlsx) to the values generated at 456a3b4.
# Applies to both the "展览" sheet (sheet 1) and the "全部类型" sheet (sheet 4),
# which carry the same exhibition rows (the latter just aggregates every category).

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$addr, [string]$text)
    # Plain Value assignment auto-parses strings that look like dates
    # (e.g. "2024-08-17") into date serials. Force a Text number format
    # first so the literal string is preserved, then drop the format back
    # to the sheet's normal (unstyled) look so no stray formatting lingers
    # on the cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).ClearFormats()
}

function Update-ExpoRow {
    param($ws, [int]$row, [hashtable]$vals)

    if ($vals.ContainsKey("B")) { Set-TextValue $ws "B$row" $vals["B"] }
    if ($vals.ContainsKey("D")) { Set-TextValue $ws "D$row" $vals["D"] }
    if ($vals.ContainsKey("E")) { Set-TextValue $ws "E$row" $vals["E"] }
    if ($vals.ContainsKey("F")) { $ws.Range("F$row").Value = $vals["F"] }
    if ($vals.ContainsKey("G")) { $ws.Range("G$row").Value = $vals["G"] }
    if ($vals.ContainsKey("I")) { Set-TextValue $ws "I$row" $vals["I"] }
}

# Row-by-row updates shared by both sheets (same events, different row numbers).
$rowUpdates = @(
    @{ F = 2010 },                                            # 2009 -> 2010
    @{ F = 322  },                                             # 320  -> 322
    @{ F = 2051 },                                             # 2050 -> 2051
    @{ F = 10405 },                                            # 10400 -> 10405
    @{ F = 270  },                                             # 268  -> 270
    @{ F = 399  },                                             # 398  -> 399
    @{ F = 7261 },                                             # 7254 -> 7261
    @{ F = 1108 },                                             # 1107 -> 1108
    @{ F = 690  },                                             # 688  -> 690
    @{ F = 153  },                                             # 149  -> 153
    @{
        B = "2024-08-17"
        D = "清禾路886号 尹山湖大剧院"
        E = "2024.08.17 10:00-08.18 17:00"
        F = 59
        G = 60
        I = "//i1.hdslb.com/bfs/openplatform/202405/1V6xl1Pg1715399710219.png"
    },
    @{ F = 1227 }                                              # 270 -> 1227
)

# Sheet "展览": rows 4,5,8,9,12,14,15,16,17,18,19,20
$sheet1RowMap = @(4, 5, 8, 9, 12, 14, 15, 16, 17, 18, 19, 20)
# Sheet "全部类型": rows 4,5,9,12,15,17,18,19,20,21,22,23
$sheet4RowMap = @(4, 5, 9, 12, 15, 17, 18, 19, 20, 21, 22, 23)

$wsExpo = $wb.Worksheets.Item(1)
for ($i = 0; $i -lt $rowUpdates.Count; $i++) {
    Update-ExpoRow $wsExpo $sheet1RowMap[$i] $rowUpdates[$i]
}

$wsAll = $wb.Worksheets.Item(4)
for ($i = 0; $i -lt $rowUpdates.Count; $i++) {
    Update-ExpoRow $wsAll $sheet4RowMap[$i] $rowUpdates[$i]
}
